# Auto-generated cell updates derived from the canonical OOXML diff.
# Each entry: Sheet, Row, Col, NewValue ($null => clear the cell)
$wb = $excel.ActiveWorkbook

$edits = @(
    ,@("ALC", 64, 8, 4456.25)
    ,@("ALC", 64, 9, 3600)
    ,@("ALC", 64, 10, 4970)
    ,@("ALC", 64, 11, 3600)
    ,@("ALC", 64, 12, 4970)
    ,@("ALC", 64, 13, -3352)
    ,@("ALC", 64, 14, -5466)
    ,@("ALC", 67, 8, 4456.25)
    ,@("ALC", 67, 9, 3600)
    ,@("ALC", 67, 10, 4970)
    ,@("ALC", 67, 11, 3600)
    ,@("ALC", 67, 12, 4970)
    ,@("ALC", 67, 13, -2742)
    ,@("ALC", 67, 14, -6686)
    ,@("ALC", 80, 8, 7613095)
    ,@("ALC", 80, 9, 660.6667)
    ,@("ALC", 80, 10, 14329949)
    ,@("ALC", 80, 11, 1982.0001)
    ,@("ALC", 80, 12, 42989847)
    ,@("ALC", 80, 13, -984.0001)
    ,@("ALC", 80, 14, -42991843)
    ,@("ALC", 83, 8, 7613095)
    ,@("ALC", 83, 9, 660.6667)
    ,@("ALC", 83, 10, 14329949)
    ,@("ALC", 83, 11, 5946.0003)
    ,@("ALC", 83, 12, 128969541)
    ,@("ALC", 83, 13, -954.0002999999997)
    ,@("ALC", 83, 14, -128979525)
    ,@("ALC", 100, 8, 71430100)
    ,@("ALC", 100, 9, 100000860)
    ,@("ALC", 100, 10, 3200)
    ,@("ALC", 100, 11, 100000860)
    ,@("ALC", 100, 12, 3200)
    ,@("ALC", 100, 13, -100000319)
    ,@("ALC", 100, 14, -4282)
    ,@("ALC", 106, 8, 11113284)
    ,@("ALC", 106, 9, 17545178)
    ,@("ALC", 106, 10, 3650.182)
    ,@("ALC", 106, 11, 17545178)
    ,@("ALC", 106, 12, 3650.182)
    ,@("ALC", 106, 13, -17544547)
    ,@("ALC", 106, 14, -4912.182)
    ,@("ALC", 112, 8, 1073.303)
    ,@("ALC", 112, 10, 1073.303)
    ,@("ALC", 112, 12, 3219.909000000001)
    ,@("ALC", 112, 14, -5435.909000000001)
    ,@("ALC", 137, 8, 2392.2856)
    ,@("ALC", 137, 9, 2199.8)
    ,@("ALC", 137, 10, 2499.2222)
    ,@("ALC", 137, 11, 6599.400000000001)
    ,@("ALC", 137, 12, 7497.6666)
    ,@("ALC", 137, 13, -4049.400000000001)
    ,@("ALC", 137, 14, -12597.6666)
    ,@("ALC", 138, 8, 2198.9468)
    ,@("ALC", 138, 9, 2127.6428)
    ,@("ALC", 138, 10, 2215.3115)
    ,@("ALC", 138, 11, 6382.928400000001)
    ,@("ALC", 138, 12, 6645.934499999999)
    ,@("ALC", 138, 13, -1242.928400000001)
    ,@("ALC", 138, 14, -16925.9345)
    ,@("ARM", 2, 8, 1634.2433)
    ,@("ARM", 2, 9, 1492.3)
    ,@("ARM", 2, 11, 1492.3)
    ,@("ARM", 2, 13, -1379.3)
    ,@("ARM", 32, 8, 5864.745)
    ,@("ARM", 32, 9, 5892.273)
    ,@("ARM", 32, 10, 5691.7144)
    ,@("ARM", 32, 11, 5892.273)
    ,@("ARM", 32, 12, 5691.7144)
    ,@("ARM", 32, 13, -5605.273)
    ,@("ARM", 32, 14, -6265.7144)
    ,@("ARM", 61, 8, 1980.591)
    ,@("ARM", 61, 9, 1719.6316)
    ,@("ARM", 61, 11, 1719.6316)
    ,@("ARM", 61, 13, -1507.6316)
    ,@("ARM", 76, 8, 14999)
    ,@("ARM", 76, 10, 14999)
    ,@("ARM", 76, 12, 14999)
    ,@("ARM", 76, 14, -15675)
    ,@("ARM", 79, 8, 14999)
    ,@("ARM", 79, 10, 14999)
    ,@("ARM", 79, 12, 14999)
    ,@("ARM", 79, 14, -17339)
    ,@("ARM", 102, 8, 1750)
    ,@("ARM", 102, 9, 0)
    ,@("ARM", 102, 10, 1750)
    ,@("ARM", 102, 11, 0)
    ,@("ARM", 102, 12, 1750)
    ,@("ARM", 102, 13, $null)
    ,@("ARM", 102, 14, -4994)
    ,@("ARM", 110, 8, 579.7143)
    ,@("ARM", 110, 9, 577.25)
    ,@("ARM", 110, 10, 583)
    ,@("ARM", 110, 11, 577.25)
    ,@("ARM", 110, 12, 583)
    ,@("ARM", 110, 13, 1467.75)
    ,@("ARM", 110, 14, -4673)
    ,@("ARM", 116, 8, 1634.2433)
    ,@("ARM", 116, 9, 1492.3)
    ,@("ARM", 116, 11, 1492.3)
    ,@("ARM", 116, 13, 801.7)
    ,@("ARM", 132, 8, 28999.072)
    ,@("ARM", 132, 9, 1432.9546)
    ,@("ARM", 132, 11, 4298.8638)
    ,@("ARM", 132, 13, -1768.8638)
    ,@("ARM", 136, 8, 1980.591)
    ,@("ARM", 136, 9, 1719.6316)
    ,@("ARM", 136, 11, 5158.8948)
    ,@("ARM", 136, 13, -2608.8948)
    ,@("BSM", 3, 8, 1634.2433)
    ,@("BSM", 3, 9, 1492.3)
    ,@("BSM", 3, 11, 1492.3)
    ,@("BSM", 3, 13, -1378.3)
    ,@("BSM", 105, 8, 3820.5833)
    ,@("BSM", 105, 9, 4264)
    ,@("BSM", 105, 10, 3199.8)
    ,@("BSM", 105, 11, 4264)
    ,@("BSM", 105, 12, 3199.8)
    ,@("BSM", 105, 13, -2517)
    ,@("BSM", 105, 14, -6693.8)
    ,@("CRP", 31, 8, 14045.314)
    ,@("CRP", 31, 9, 22680.053)
    ,@("CRP", 31, 11, 22680.053)
    ,@("CRP", 31, 13, -22385.053)
    ,@("CRP", 34, 8, 14045.314)
    ,@("CRP", 34, 9, 22680.053)
    ,@("CRP", 34, 11, 22680.053)
    ,@("CRP", 34, 13, -22478.053)
    ,@("CRP", 58, 8, 27792.895)
    ,@("CRP", 58, 9, 1634.6923)
    ,@("CRP", 58, 10, 84469)
    ,@("CRP", 58, 11, 1634.6923)
    ,@("CRP", 58, 12, 84469)
    ,@("CRP", 58, 13, -1431.6923)
    ,@("CRP", 58, 14, -84875)
    ,@("CRP", 105, 8, 11364423)
    ,@("CRP", 105, 9, 11364423)
    ,@("CRP", 105, 10, 0)
    ,@("CRP", 105, 11, 11364423)
    ,@("CRP", 105, 12, 0)
    ,@("CRP", 105, 13, -11362676)
    ,@("CRP", 105, 14, $null)
    ,@("CRP", 132, 8, 20682.5)
    ,@("CRP", 132, 9, 22670.166)
    ,@("CRP", 132, 11, 68010.49800000001)
    ,@("CRP", 132, 13, -65480.49800000001)
    ,@("CRP", 136, 8, 27792.895)
    ,@("CRP", 136, 9, 1634.6923)
    ,@("CRP", 136, 10, 84469)
    ,@("CRP", 136, 11, 4904.0769)
    ,@("CRP", 136, 12, 253407)
    ,@("CRP", 136, 13, -2354.0769)
    ,@("CRP", 136, 14, -258507)
    ,@("CUL", 113, 8, 658.82355)
    ,@("CUL", 113, 9, 437.77777)
    ,@("CUL", 113, 10, 907.5)
    ,@("CUL", 113, 11, 1313.33331)
    ,@("CUL", 113, 12, 2722.5)
    ,@("CUL", 113, 13, 856.66669)
    ,@("CUL", 113, 14, -7062.5)
    ,@("CUL", 118, 8, 50002210)
    ,@("CUL", 118, 10, 4199.8)
    ,@("CUL", 118, 12, 12599.4)
    ,@("CUL", 118, 14, -15085.4)
    ,@("CUL", 121, 8, 5212036.5)
    ,@("CUL", 121, 9, 399.875)
    ,@("CUL", 121, 10, 6949249)
    ,@("CUL", 121, 11, 1199.625)
    ,@("CUL", 121, 12, 20847747)
    ,@("CUL", 121, 13, 110.375)
    ,@("CUL", 121, 14, -20850367)
    ,@("CUL", 122, 8, 708.5)
    ,@("CUL", 122, 9, 368)
    ,@("CUL", 122, 10, 859.8333)
    ,@("CUL", 122, 11, 3312)
    ,@("CUL", 122, 12, 7738.4997)
    ,@("CUL", 122, 13, -862)
    ,@("CUL", 122, 14, -12638.4997)
    ,@("CUL", 132, 8, 1158.9)
    ,@("CUL", 132, 9, 500)
    ,@("CUL", 132, 11, 4500)
    ,@("CUL", 132, 13, -1970)
    ,@("GSM", 102, 8, 31253150)
    ,@("GSM", 102, 9, 33336528)
    ,@("GSM", 102, 10, 2500)
    ,@("GSM", 102, 11, 33336528)
    ,@("GSM", 102, 12, 2500)
    ,@("GSM", 102, 13, -33334906)
    ,@("GSM", 102, 14, -5744)
    ,@("GSM", 132, 8, 19371.162)
    ,@("GSM", 132, 9, 3701.5)
    ,@("GSM", 132, 10, 47861.453)
    ,@("GSM", 132, 11, 11104.5)
    ,@("GSM", 132, 12, 143584.359)
    ,@("GSM", 132, 13, -8574.5)
    ,@("GSM", 132, 14, -148644.359)
    ,@("LTW", 68, 8, 2313)
    ,@("LTW", 68, 9, 2119.9)
    ,@("LTW", 68, 11, 2119.9)
    ,@("LTW", 68, 13, -1370.9)
    ,@("LTW", 71, 8, 2313)
    ,@("LTW", 71, 9, 2119.9)
    ,@("LTW", 71, 11, 10599.5)
    ,@("LTW", 71, 13, -6855.5)
    ,@("LTW", 122, 8, 786825.4)
    ,@("LTW", 122, 9, 1402809.6)
    ,@("LTW", 122, 10, 2845.4546)
    ,@("LTW", 122, 11, 4208428.800000001)
    ,@("LTW", 122, 12, 8536.363799999999)
    ,@("LTW", 122, 13, -4205978.800000001)
    ,@("LTW", 122, 14, -13436.3638)
    ,@("LTW", 132, 8, 2592.111)
    ,@("LTW", 132, 9, 2373.4614)
    ,@("LTW", 132, 11, 7120.3842)
    ,@("LTW", 132, 13, -4590.3842)
    ,@("LTW", 136, 8, 37682)
    ,@("LTW", 136, 9, 63929.75)
    ,@("LTW", 136, 10, 2685)
    ,@("LTW", 136, 11, 191789.25)
    ,@("LTW", 136, 12, 8055)
    ,@("LTW", 136, 13, -189239.25)
    ,@("LTW", 136, 14, -13155)
    ,@("WVR", 14, 8, 102490)
    ,@("WVR", 14, 10, 4980)
    ,@("WVR", 14, 12, 4980)
    ,@("WVR", 14, 14, -5316)
    ,@("WVR", 62, 8, 3439)
    ,@("WVR", 62, 9, 3133.5557)
    ,@("WVR", 62, 10, 3744.4443)
    ,@("WVR", 62, 11, 3133.5557)
    ,@("WVR", 62, 12, 3744.4443)
    ,@("WVR", 62, 13, -2509.5557)
    ,@("WVR", 62, 14, -4992.4443)
    ,@("WVR", 65, 8, 3439)
    ,@("WVR", 65, 9, 3133.5557)
    ,@("WVR", 65, 10, 3744.4443)
    ,@("WVR", 65, 11, 15667.7785)
    ,@("WVR", 65, 12, 18722.2215)
    ,@("WVR", 65, 13, -12547.7785)
    ,@("WVR", 65, 14, -24962.2215)
    ,@("WVR", 96, 8, 999)
    ,@("WVR", 96, 10, 1038.8)
    ,@("WVR", 96, 12, 1038.8)
    ,@("WVR", 96, 14, -3784.8)
    ,@("WVR", 100, 8, 487.5)
    ,@("WVR", 100, 9, 491.66666)
    ,@("WVR", 100, 10, 475)
    ,@("WVR", 100, 11, 983.33332)
    ,@("WVR", 100, 12, 950)
    ,@("WVR", 100, 13, -442.33332)
    ,@("WVR", 100, 14, -2032)
    ,@("WVR", 132, 8, 952.6177)
    ,@("WVR", 132, 9, 649.6539)
    ,@("WVR", 132, 11, 1948.9617)
    ,@("WVR", 136, 8, 41668664)
    ,@("WVR", 136, 9, 50001852)
    ,@("WVR", 136, 10, 2725)
    ,@("WVR", 136, 11, 150005556)
    ,@("WVR", 136, 12, 8175)
    ,@("WVR", 136, 13, -150003006)
    ,@("WVR", 136, 14, -13275)
)

foreach ($edit in $edits) {
    $sheetName = $edit[0]
    $row = $edit[1]
    $col = $edit[2]
    $newValue = $edit[3]
    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Cells.Item($row, $col)
    if ($null -eq $newValue) {
        $cell.ClearContents()
    } else {
        $cell.Value = $newValue
    }
}

Write-Host "Applied $($edits.Count) cell edits."